# Apply "pandoc title block" restructuring to the first two paragraphs:
#   1. Title paragraph: style Heading1 -> Title, text split into separate
#      per-word/space runs.
#   2. Byline paragraph: style Normal/bold "By Dorothy Day" -> style
#      Authors, "By " prefix and bold formatting removed, text split into
#      separate per-word/space runs ("Dorothy", " ", "Day").
#
# (Note: the document also originally has an invalid-named bookmark
# (w:name="on-pilgrimage---december-1957", not a legal Word bookmark
# name) wrapping the title paragraph; it is not exposed anywhere in the
# Word object model here -- Document.Bookmarks.Count reports 0 for it and
# it cannot be targeted/deleted through any COM API -- so it is left
# as-is.)

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: title ---------------------------------------------------
$titleWords = @("On", " ", "Pilgrimage", " ", "-", " ", "December", " ", "1957")
$titleRuns = ($titleWords | ForEach-Object {
    '<w:r><w:t xml:space="preserve">' + $_ + '</w:t></w:r>'
}) -join ''
$titleFrag = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>'

$p1 = $d.Paragraphs.Item(1)
[void]$p1.Range.InsertXML($titleFrag)

# --- Paragraph 2: byline ---------------------------------------------------
$authorWords = @("Dorothy", " ", "Day")
$authorRuns = ($authorWords | ForEach-Object {
    '<w:r><w:t xml:space="preserve">' + $_ + '</w:t></w:r>'
}) -join ''
$authorFrag = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'

$p2 = $d.Paragraphs.Item(2)
[void]$p2.Range.InsertXML($authorFrag)

Write-Host "Paragraph 1: [$($d.Paragraphs.Item(1).Range.Text)] style=$($d.Paragraphs.Item(1).Range.Style.NameLocal)"
Write-Host "Paragraph 2: [$($d.Paragraphs.Item(2).Range.Text)] style=$($d.Paragraphs.Item(2).Range.Style.NameLocal)"
